# Auto-generated Excel COM-interop script
# Applies updated market-price / profit figures to the per-job "Profits" sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) as captured by the scheduled price-update run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2125
$ws.Range("I17").Value = 2125
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 6375
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -6207
$ws.Range("N17").ClearContents()

$ws.Range("H33").Value = 838.46155
$ws.Range("I33").Value = 939.2727
$ws.Range("K33").Value = 939.2727
$ws.Range("M33").Value = -710.2727

$ws.Range("H112").Value = 1171.1482
$ws.Range("J112").Value = 1158.5
$ws.Range("L112").Value = 3475.5
$ws.Range("N112").Value = -5691.5

$ws.Range("H133").Value = 85999
$ws.Range("J133").Value = 85999
$ws.Range("L133").Value = 85999
$ws.Range("N133").Value = -96119

$ws.Range("H138").Value = 3588.75
$ws.Range("I138").Value = 2880.8572
$ws.Range("J138").Value = 3864.0417
$ws.Range("K138").Value = 8642.571599999999
$ws.Range("L138").Value = 11592.1251
$ws.Range("M138").Value = -3502.571599999999
$ws.Range("N138").Value = -21872.1251

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4923.457
$ws.Range("I32").Value = 4163.5
$ws.Range("J32").Value = 13029.667
$ws.Range("K32").Value = 4163.5
$ws.Range("L32").Value = 13029.667
$ws.Range("M32").Value = -3876.5
$ws.Range("N32").Value = -13603.667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3491.7058
$ws.Range("I99").Value = 2784.9
$ws.Range("K99").Value = 2784.9
$ws.Range("M99").Value = -1286.9

$ws.Range("H105").Value = 2911.3333
$ws.Range("I105").Value = 2774.0588
$ws.Range("K105").Value = 2774.0588
$ws.Range("M105").Value = -1027.0588

$ws.Range("H134").Value = 1636.7059
$ws.Range("I134").Value = 1813.3572
$ws.Range("K134").Value = 5440.071599999999
$ws.Range("M134").Value = -2905.071599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 9999999
$ws.Range("I4").Value = 9999999
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 9999999
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -9999887
$ws.Range("N4").ClearContents()

$ws.Range("H7").Value = 68.13333
$ws.Range("I7").Value = 64.333336
$ws.Range("J7").Value = 83.333336
$ws.Range("K7").Value = 64.333336
$ws.Range("L7").Value = 83.333336
$ws.Range("M7").Value = 48.666664
$ws.Range("N7").Value = -309.333336

$ws.Range("H22").Value = 9468.444
$ws.Range("I22").Value = 210.75
$ws.Range("J22").Value = 35919
$ws.Range("K22").Value = 210.75
$ws.Range("L22").Value = 35919
$ws.Range("M22").Value = 139.25
$ws.Range("N22").Value = -36619

$ws.Range("H33").Value = 2174.8333
$ws.Range("I33").Value = 2174.8333
$ws.Range("K33").Value = 2174.8333
$ws.Range("M33").Value = -1795.8333

$ws.Range("H134").Value = 1835.4722
$ws.Range("I134").Value = 1583.7812
$ws.Range("K134").Value = 4751.3436
$ws.Range("M134").Value = -2216.3436

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 11040759
$ws.Range("I4").Value = 11040759
$ws.Range("K4").Value = 33122277
$ws.Range("M4").Value = -33122165

$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()

$ws.Range("H127").Value = 1998.6666
$ws.Range("J127").Value = 1998.6666
$ws.Range("L127").Value = 5995.9998
$ws.Range("N127").Value = -15915.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 827.6667
$ws.Range("I6").Value = 827.6667
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 827.6667
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -714.6667
$ws.Range("N6").ClearContents()

$ws.Range("H16").Value = 827.6667
$ws.Range("I16").Value = 827.6667
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 827.6667
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -577.6667
$ws.Range("N16").ClearContents()

$ws.Range("H43").Value = 2999.5
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 2999.5
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 2999.5
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -3301.5

$ws.Range("H70").Value = 5299.5625
$ws.Range("J70").Value = 5756.875
$ws.Range("L70").Value = 5756.875
$ws.Range("N70").Value = -6296.875

$ws.Range("H73").Value = 5299.5625
$ws.Range("J73").Value = 5756.875
$ws.Range("L73").Value = 5756.875
$ws.Range("N73").Value = -7628.875

$ws.Range("H80").Value = 12548.728
$ws.Range("I80").Value = 3333.3333
$ws.Range("J80").Value = 16004.5
$ws.Range("K80").Value = 3333.3333
$ws.Range("L80").Value = 16004.5
$ws.Range("M80").Value = -2335.3333
$ws.Range("N80").Value = -18000.5

$ws.Range("H83").Value = 12548.728
$ws.Range("I83").Value = 3333.3333
$ws.Range("J83").Value = 16004.5
$ws.Range("K83").Value = 16666.6665
$ws.Range("L83").Value = 80022.5
$ws.Range("M83").Value = -11674.6665
$ws.Range("N83").Value = -90006.5

$ws.Range("H92").Value = 26041
$ws.Range("J92").Value = 26041
$ws.Range("L92").Value = 26041
$ws.Range("N92").Value = -29785

$ws.Range("H113").Value = 4000
$ws.Range("I113").Value = 4000
$ws.Range("K113").Value = 4000
$ws.Range("M113").Value = -1830

$ws.Range("H134").Value = 30000
$ws.Range("J134").Value = 30000
$ws.Range("L134").Value = 90000
$ws.Range("N134").Value = -95070

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4193.5
$ws.Range("I132").Value = 4103.8887
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 12311.6661
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -9781.666100000002
$ws.Range("N132").Value = -20060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 16000000
$ws.Range("I2").Value = 31000000
$ws.Range("K2").Value = 31000000
$ws.Range("M2").Value = -30999888

$ws.Range("H9").Value = 1500
$ws.Range("J9").Value = 2000
$ws.Range("L9").Value = 2000
$ws.Range("N9").Value = -2280

$ws.Range("H12").Value = 1004
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 1004
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 1004
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -1288

$ws.Range("H34").Value = 39999
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()

$ws.Range("H61").Value = 38498.5
$ws.Range("J61").Value = 38498.332
$ws.Range("L61").Value = 38498.332
$ws.Range("N61").Value = -39082.332

$ws.Range("H100").Value = 6251038
$ws.Range("I100").Value = 9091778
$ws.Range("J100").Value = 1409.8
$ws.Range("K100").Value = 18183556
$ws.Range("L100").Value = 2819.6
$ws.Range("M100").Value = -18183015
$ws.Range("N100").Value = -3901.6

$ws.Range("H122").Value = 2350.7778
$ws.Range("I122").Value = 2457.4375
$ws.Range("K122").Value = 7372.3125
$ws.Range("M122").Value = -4922.3125

$ws.Range("H126").Value = 4672.7896
$ws.Range("I126").Value = 2426.5454
$ws.Range("K126").Value = 7279.6362
$ws.Range("M126").Value = -4809.6362

$ws.Range("H132").Value = 13330
$ws.Range("I132").Value = 9990
$ws.Range("K132").Value = 29970
$ws.Range("M132").Value = -27440

$ws.Range("H136").Value = 1558.6666
$ws.Range("I136").Value = 1558.6666
$ws.Range("K136").Value = 4675.9998
$ws.Range("M136").Value = -2125.9998
